# Update "想去人数" (F) and "最低票价" (G) figures in both the
# "展览" (Exhibition) and "全部类型" (All types) sheets, since they contain
# duplicate listings of the same events.

$wb = $excel.ActiveWorkbook

# Row -> [F value, G value (or $null to leave unchanged)]
$sheet1Updates = @{
    2  = @(623, $null)
    4  = @(1298, $null)
    6  = @(14345, $null)
    7  = @(16569, 65)
    9  = @(98, $null)
    10 = @(13, $null)
    19 = @(105, $null)
    23 = @(70, $null)
    24 = @(40, $null)
    25 = @(14, $null)
    27 = @(6741, $null)
    30 = @(1120, $null)
    32 = @(9, $null)
    33 = @(5753, $null)
    35 = @(146, $null)
    36 = @(191, $null)
    37 = @(4830, $null)
}

$sheet4Updates = @{
    2  = @(623, $null)
    4  = @(1298, $null)
    6  = @(14345, $null)
    7  = @(16569, 65)
    9  = @(98, $null)
    10 = @(13, $null)
    19 = @(105, $null)
    23 = @(70, $null)
    25 = @(40, $null)
    26 = @(14, $null)
    28 = @(6741, $null)
    31 = @(1120, $null)
    33 = @(9, $null)
    36 = @(5753, $null)
    38 = @(146, $null)
    39 = @(191, $null)
    40 = @(4830, $null)
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $vals = $sheet1Updates[$row]
    $ws1.Range("F$row").Value = $vals[0]
    if ($null -ne $vals[1]) {
        $ws1.Range("G$row").Value = $vals[1]
    }
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $vals = $sheet4Updates[$row]
    $ws4.Range("F$row").Value = $vals[0]
    if ($null -ne $vals[1]) {
        $ws4.Range("G$row").Value = $vals[1]
    }
}
